# Updated TPM-derived expression values for the Inhba-Bambi ligand-receptor pair sheet.
# Only the "ECs" cluster's raw Ligand (Inhba) and Receptor (Bambi) expression values
# changed (new TPM normalization); every other column on each row is a value derived
# from those raw numbers (cluster-specificity ratios and ligand*receptor edge weights).
# The literal replacement values below are the recomputed results for each affected cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.08097566666666667  # Ligand average expression value: 3.675031333333333 -> 0.08097566666666667
$ws.Range("H2").Value = 0.242927  # Ligand total expression value: 11.025094 -> 0.242927
$ws.Range("I2").Value = 0.005588990034505014  # Ligand derived specificity of average expression value: 0.2032371147293133 -> 0.005588990034505014
$ws.Range("J2").Value = 0.005588990034505015  # Ligand derived specificity of total expression value: 0.2032371147293133 -> 0.005588990034505015
$ws.Range("M2").Value = 2.402487333333333  # Receptor average expression value: 2.035948 -> 2.402487333333333
$ws.Range("N2").Value = 7.207462  # Receptor total expression value: 6.107844 -> 7.207462
$ws.Range("O2").Value = 0.5333179916503753  # Receptor derived specificity of average expression value: 0.4919823674428878 -> 0.5333179916503753
$ws.Range("P2").Value = 0.5333179916503754  # Receptor derived specificity of total expression value: 0.4919823674428878 -> 0.5333179916503754
$ws.Range("Q2").Value = 0.1945430134748889  # Edge average expression weight: 7.482172693037333 -> 0.1945430134748889
$ws.Range("R2").Value = 1.750887121274  # Edge total expression weight: 67.339554237336 -> 1.750887121274
$ws.Range("S2").Value = 0.002980708940556176  # Edge average expression derived specificity: 0.09998907685678936 -> 0.002980708940556176
$ws.Range("T2").Value = 0.002980708940556177  # Edge total expression derived specificity: 0.09998907685678936 -> 0.002980708940556177

# Row 3
$ws.Range("G3").Value = 0.08097566666666667  # Ligand average expression value: 3.675031333333333 -> 0.08097566666666667
$ws.Range("H3").Value = 0.242927  # Ligand total expression value: 11.025094 -> 0.242927
$ws.Range("I3").Value = 0.005588990034505014  # Ligand derived specificity of average expression value: 0.2032371147293133 -> 0.005588990034505014
$ws.Range("J3").Value = 0.005588990034505015  # Ligand derived specificity of total expression value: 0.2032371147293133 -> 0.005588990034505015
$ws.Range("O3").Value = 0.2964836714669855  # Receptor derived specificity of average expression value: 0.3227442459227168 -> 0.2964836714669855
$ws.Range("P3").Value = 0.2964836714669856  # Receptor derived specificity of total expression value: 0.3227442459227168 -> 0.2964836714669856
$ws.Range("Q3").Value = 0.1081509114567778  # Edge average expression weight: 4.908363273726889 -> 0.1081509114567778
$ws.Range("R3").Value = 0.9733582031110001  # Edge total expression weight: 44.17526946354199 -> 0.9733582031110001
$ws.Range("S3").Value = 0.001657044285222441  # Edge average expression derived specificity: 0.0655936093368209 -> 0.001657044285222441
$ws.Range("T3").Value = 0.001657044285222441  # Edge total expression derived specificity: 0.0655936093368209 -> 0.001657044285222441

# Row 4
$ws.Range("G4").Value = 0.08097566666666667  # Ligand average expression value: 3.675031333333333 -> 0.08097566666666667
$ws.Range("H4").Value = 0.242927  # Ligand total expression value: 11.025094 -> 0.242927
$ws.Range("I4").Value = 0.005588990034505014  # Ligand derived specificity of average expression value: 0.2032371147293133 -> 0.005588990034505014
$ws.Range("J4").Value = 0.005588990034505015  # Ligand derived specificity of total expression value: 0.2032371147293133 -> 0.005588990034505015
$ws.Range("M4").Value = 0.7667083333333334  # Receptor average expression value: 0.7667083333333333 -> 0.7667083333333334
$ws.Range("O4").Value = 0.1701983368826391  # Receptor derived specificity of average expression value: 0.1852733866343954 -> 0.1701983368826391
$ws.Range("P4").Value = 0.1701983368826391  # Receptor derived specificity of total expression value: 0.1852733866343954 -> 0.1701983368826391
$ws.Range("Q4").Value = 0.06208471843055557  # Edge average expression weight: 2.817677148527777 -> 0.06208471843055557
$ws.Range("R4").Value = 0.5587624658750001  # Edge total expression weight: 25.35909433675 -> 0.5587624658750001
$ws.Range("S4").Value = 0.0009512368087263971  # Edge average expression derived specificity: 0.03765442853570304 -> 0.0009512368087263971
$ws.Range("T4").Value = 0.0009512368087263974  # Edge total expression derived specificity: 0.03765442853570304 -> 0.0009512368087263974

# Row 5
$ws.Range("I5").Value = 0.6976944377922635  # Ligand derived specificity of average expression value: 0.5590213983169419 -> 0.6976944377922635
$ws.Range("J5").Value = 0.6976944377922635  # Ligand derived specificity of total expression value: 0.5590213983169419 -> 0.6976944377922635
$ws.Range("M5").Value = 2.402487333333333  # Receptor average expression value: 2.035948 -> 2.402487333333333
$ws.Range("N5").Value = 7.207462  # Receptor total expression value: 6.107844 -> 7.207462
$ws.Range("O5").Value = 0.5333179916503753  # Receptor derived specificity of average expression value: 0.4919823674428878 -> 0.5333179916503753
$ws.Range("P5").Value = 0.5333179916503754  # Receptor derived specificity of total expression value: 0.4919823674428878 -> 0.5333179916503754
$ws.Range("Q5").Value = 24.285528794076  # Edge average expression weight: 20.580368142312 -> 24.285528794076
$ws.Range("R5").Value = 218.569759146684  # Edge total expression weight: 185.223313280808 -> 218.569759146684
$ws.Range("S5").Value = 0.3720929963490077  # Edge average expression derived specificity: 0.2750286709952026 -> 0.3720929963490077
$ws.Range("T5").Value = 0.3720929963490077  # Edge total expression derived specificity: 0.2750286709952026 -> 0.3720929963490077

# Row 6
$ws.Range("I6").Value = 0.6976944377922635  # Ligand derived specificity of average expression value: 0.5590213983169419 -> 0.6976944377922635
$ws.Range("J6").Value = 0.6976944377922635  # Ligand derived specificity of total expression value: 0.5590213983169419 -> 0.6976944377922635
$ws.Range("O6").Value = 0.2964836714669855  # Receptor derived specificity of average expression value: 0.3227442459227168 -> 0.2964836714669855
$ws.Range("P6").Value = 0.2964836714669856  # Receptor derived specificity of total expression value: 0.3227442459227168 -> 0.2964836714669856
$ws.Range("S6").Value = 0.2068550084787446  # Edge average expression derived specificity: 0.1804209396544641 -> 0.2068550084787446
$ws.Range("T6").Value = 0.2068550084787446  # Edge total expression derived specificity: 0.1804209396544641 -> 0.2068550084787446

# Row 7
$ws.Range("I7").Value = 0.6976944377922635  # Ligand derived specificity of average expression value: 0.5590213983169419 -> 0.6976944377922635
$ws.Range("J7").Value = 0.6976944377922635  # Ligand derived specificity of total expression value: 0.5590213983169419 -> 0.6976944377922635
$ws.Range("M7").Value = 0.7667083333333334  # Receptor average expression value: 0.7667083333333333 -> 0.7667083333333334
$ws.Range("O7").Value = 0.1701983368826391  # Receptor derived specificity of average expression value: 0.1852733866343954 -> 0.1701983368826391
$ws.Range("P7").Value = 0.1701983368826391  # Receptor derived specificity of total expression value: 0.1852733866343954 -> 0.1701983368826391
$ws.Range("R7").Value = 69.75239928525002  # Edge total expression weight: 69.75239928525001 -> 69.75239928525002
$ws.Range("S7").Value = 0.1187464329645112  # Edge average expression derived specificity: 0.1035717876672751 -> 0.1187464329645112
$ws.Range("T7").Value = 0.1187464329645112  # Edge total expression derived specificity: 0.1035717876672751 -> 0.1187464329645112

# Row 8
$ws.Range("I8").Value = 0.2967165721732315  # Ligand derived specificity of average expression value: 0.2377414869537448 -> 0.2967165721732315
$ws.Range("J8").Value = 0.2967165721732316  # Ligand derived specificity of total expression value: 0.2377414869537448 -> 0.2967165721732316
$ws.Range("M8").Value = 2.402487333333333  # Receptor average expression value: 2.035948 -> 2.402487333333333
$ws.Range("N8").Value = 7.207462  # Receptor total expression value: 6.107844 -> 7.207462
$ws.Range("O8").Value = 0.5333179916503753  # Receptor derived specificity of average expression value: 0.4919823674428878 -> 0.5333179916503753
$ws.Range("P8").Value = 0.5333179916503754  # Receptor derived specificity of total expression value: 0.4919823674428878 -> 0.5333179916503754
$ws.Range("Q8").Value = 10.32818733655733  # Edge average expression weight: 8.752450870288 -> 10.32818733655733
$ws.Range("R8").Value = 92.953686029016  # Edge total expression weight: 78.77205783259201 -> 92.953686029016
$ws.Range("S8").Value = 0.1582442863608114  # Edge average expression derived specificity: 0.1169646195908958 -> 0.1582442863608114
$ws.Range("T8").Value = 0.1582442863608115  # Edge total expression derived specificity: 0.1169646195908958 -> 0.1582442863608115

# Row 9
$ws.Range("I9").Value = 0.2967165721732315  # Ligand derived specificity of average expression value: 0.2377414869537448 -> 0.2967165721732315
$ws.Range("J9").Value = 0.2967165721732316  # Ligand derived specificity of total expression value: 0.2377414869537448 -> 0.2967165721732316
$ws.Range("O9").Value = 0.2964836714669855  # Receptor derived specificity of average expression value: 0.3227442459227168 -> 0.2964836714669855
$ws.Range("P9").Value = 0.2964836714669856  # Receptor derived specificity of total expression value: 0.3227442459227168 -> 0.2964836714669856
$ws.Range("S9").Value = 0.08797161870301846  # Edge average expression derived specificity: 0.07672969693143178 -> 0.08797161870301846
$ws.Range("T9").Value = 0.08797161870301849  # Edge total expression derived specificity: 0.07672969693143178 -> 0.08797161870301849

# Row 10
$ws.Range("I10").Value = 0.2967165721732315  # Ligand derived specificity of average expression value: 0.2377414869537448 -> 0.2967165721732315
$ws.Range("J10").Value = 0.2967165721732316  # Ligand derived specificity of total expression value: 0.2377414869537448 -> 0.2967165721732316
$ws.Range("M10").Value = 0.7667083333333334  # Receptor average expression value: 0.7667083333333333 -> 0.7667083333333334
$ws.Range("O10").Value = 0.1701983368826391  # Receptor derived specificity of average expression value: 0.1852733866343954 -> 0.1701983368826391
$ws.Range("P10").Value = 0.1701983368826391  # Receptor derived specificity of total expression value: 0.1852733866343954 -> 0.1701983368826391
$ws.Range("R10").Value = 29.66440850850001  # Edge total expression weight: 29.6644085085 -> 29.66440850850001
$ws.Range("S10").Value = 0.05050066710940156  # Edge average expression derived specificity: 0.04404717043141723 -> 0.05050066710940156
$ws.Range("T10").Value = 0.05050066710940157  # Edge total expression derived specificity: 0.04404717043141723 -> 0.05050066710940157
